$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 4 new data rows (8..11) below the existing template row 7,
#        pushing the "total" row (old 8) and "footer" row (old 9) down to 12/13.
$ws.Rows("8:11").Insert()

# --- 2. Copy row 7's formatting (fonts/fills/borders/number formats) onto
#        the freshly inserted rows so they reuse the same cell styles.
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)
$ws.Range("A7:Q7").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)
$ws.Range("A7:Q7").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)
$ws.Range("A7:Q7").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Re-create the merged cells for rows 8..11 (same pattern as row 7).
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

# --- 4. Row heights (match the captured diff; row 8/10 keep the default
#        "24.75" template height, 7/9/11 use "25.5").
$ws.Rows("7").RowHeight = 25.5
$ws.Rows("8").RowHeight = 24.75
$ws.Rows("9").RowHeight = 25.5
$ws.Rows("10").RowHeight = 24.75
$ws.Rows("11").RowHeight = 25.5

# --- 5. Fill in the five missing-items data rows.
# Row 7 : DICLOSP S.R. 75MG 30 F.C. TABS
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "DICLOSP S.R. 75MG 30 F.C. TABS"
$ws.Range("H7").Value = "0:0"
$ws.Range("L7").Value = "1"
$ws.Range("N7").Value = "108.00"
$ws.Range("P7").Value = "35.6400"
$ws.Range("Q7").Value = "0:1"

# Row 8 : LIPONA 20MG 10 F.C.TAB.
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "LIPONA 20MG 10 F.C.TAB."
$ws.Range("H8").Value = "0:0"
$ws.Range("L8").Value = "1"
$ws.Range("N8").Value = "41.00"
$ws.Range("P8").Value = "41.0000"
$ws.Range("Q8").Value = "1:0"

# Row 9 : STRINGAZOLE 40MG 21 ENTERIC COATED TABLETS
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "STRINGAZOLE 40MG 21 ENTERIC COATED TABLETS"
$ws.Range("H9").Value = "1:0"
$ws.Range("L9").Value = "1"
$ws.Range("N9").Value = "126.00"
$ws.Range("P9").Value = "41.5800"
$ws.Range("Q9").Value = "0:1"

# Row 10 : سرنجات انسولين
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "سرنجات انسولين"
$ws.Range("H10").Value = "9:0"
$ws.Range("L10").Value = "0"
$ws.Range("N10").Value = "7.00"
$ws.Range("P10").Value = "7.0000"
$ws.Range("Q10").Value = "1:0"

# Row 11 : مسك الرمان
$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "مسك الرمان"
$ws.Range("H11").Value = "19:0"
$ws.Range("L11").Value = "0"
$ws.Range("N11").Value = "30.00"
$ws.Range("P11").Value = "30.0000"
$ws.Range("Q11").Value = "1:0"

# --- 6. Mark the text-like numeric columns as Text format (numFmtId 49)
#        the way the source workbook stores them.
$ws.Range("C7:G11").NumberFormat = "@"
$ws.Range("H7:K11").NumberFormat = "@"
$ws.Range("N7:O11").NumberFormat = "@"
$ws.Range("Q7:Q11").NumberFormat = "@"

# --- 7. Totals row (now row 12) and the timestamp/footer row (now row 13).
$ws.Range("P12").Value = 155.22

$ws.Range("A13").Value = "Thursday, 11 September, 2025 10:48 AM"
